# The edit re-sorts the weekly price rows (rows 2-17, columns A:T) of the
# sheet into a new order. Columns A,B,C,E,F,G,H,I,J,K,L are identical on
# every data row, so the visible effect is that columns D and M:T move to
# a new row according to the mapping below (new row -> source row, using
# the ORIGINAL/before-edit row numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 17

# newRow -> oldRow (both refer to the row numbers as they exist BEFORE
# this script runs, i.e. the original layout in before.xlsx)
$mapping = @{
    2  = 10
    3  = 17
    4  = 12
    5  = 15
    6  = 7
    7  = 3
    8  = 13
    9  = 14
    10 = 16
    11 = 4
    12 = 2
    13 = 5
    14 = 8
    15 = 6
    16 = 11
    17 = 9
}

# Snapshot every source row (A:T) BEFORE writing anything, so that writes
# to one row never clobber data still needed for another row.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $snapshot[$r] = $ws.Range("A$r`:T$r").Value2
}

for ($newR = $firstRow; $newR -le $lastRow; $newR++) {
    $oldR = $mapping[$newR]
    $ws.Range("A$newR`:T$newR").Value2 = $snapshot[$oldR]
}
